$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds date-like text ("2025-09-11"); a plain Value assignment
# would have Excel auto-convert it to a date serial. Prefix with an
# apostrophe so it is stored as literal text, matching the existing rows.
$ws.Range("A18").Value = "'2025-09-11"
$ws.Range("B18").Value = "21:19:58"
$ws.Range("C18").Value = "1.00 EUR = 1667.5922 ARS"
